# Apply the commit's edit to the last paragraph of the document:
#   " в активном режиме, необходимо 3 " -> " в активном режиме, необходимо как минимум 5 "
#   " адреса."                           -> " адресов, при входящем трафике 1.5 Гбит/с и
#                                              физических линках хостов в 1 Гбит/с. Чтобы
#                                              обеспечить достаточную пропускную способность."
# with the (hidden) "_GoBack" bookmark ending up right before the final
# "пропускную способность." sentence, same as in the source paragraph.

$d = $word.ActiveDocument

# --- 1) "необходимо 3 " -> "необходимо как минимум 5 " ------------------------------
# This phrase is unique in the document, and it sits entirely inside one run
# (the run just before the separately-formatted "VIP" run), so a plain
# Find/Replace keeps that run's formatting and leaves "VIP" untouched.
$rng1 = $d.Content
$found1 = $rng1.Find.Execute( `
    " в активном режиме, необходимо 3 ", $true, $false, $false, $false, $false, `
    $true, 1, $false, " в активном режиме, необходимо как минимум 5 ", 2)

# --- 2) " адреса." -> extended sentence --------------------------------------------
# " адреса." is NOT unique in the document (it also occurs much earlier), so first
# locate the first occurrence and then search again in the remainder of the story to
# reliably land on the one that follows "VIP" at the end of the document.
$probe = $d.Content
$null = $probe.Find.Execute(" адреса.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$tailStory = $d.Range($probe.End, $d.Content.End)

$newTail = " адресов, при входящем трафике 1.5 Гбит/с и физических линках хостов в 1 Гбит/с. Чтобы обеспечить достаточную пропускную способность."
$found2 = $tailStory.Find.Execute( `
    " адреса.", $true, $false, $false, $false, $false, $true, 1, $false, $newTail, 1)

# --- 3) Re-insert the "_GoBack" bookmark right before "пропускную способность." ----
# The plain Find/Replace above folds the whole new tail into a single run, which
# pushes the paragraph's existing "_GoBack" bookmark to the very end of the
# paragraph. Re-anchor it between "... достаточную " and "пропускную способность."
# (its original relative position), which also cleanly splits that run in two so
# both halves keep the inherited "Times New Roman"/sz24 run formatting.
$anchor = $d.Content
$foundAnchor = $anchor.Find.Execute("пропускную способность.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($foundAnchor) {
    $d.Bookmarks.Add("_GoBack", $d.Range($anchor.Start, $anchor.Start))
}

Write-Output "part1=$found1 part2=$found2 bookmark=$foundAnchor"
